# Populate Cocktail Sort data (rows 37-41) on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Trial 1 (row 37)
$ws.Range("P37").Value = 2
$ws.Range("Q37").Value = 8
$ws.Range("R37").Value = 126
$ws.Range("S37").Value = 13569
$ws.Range("T37").Value = "Unmeasureable"
$ws.Range("U37").Value = "Unmeasureable"

# Trial 2 (row 38)
$ws.Range("P38").Value = 1
$ws.Range("Q38").Value = 3
$ws.Range("R38").Value = 109
$ws.Range("S38").Value = 13801
$ws.Range("T38").Value = "Unmeasureable"
$ws.Range("U38").Value = "Unmeasureable"

# Trial 3 (row 39)
$ws.Range("P39").Value = 1
$ws.Range("Q39").Value = 2
$ws.Range("R39").Value = 109
$ws.Range("S39").Value = 12396
$ws.Range("T39").Value = "Unmeasureable"
$ws.Range("U39").Value = "Unmeasureable"

# Trial 4 (row 40)
$ws.Range("P40").Value = 1
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = 113
$ws.Range("S40").Value = 12261
$ws.Range("T40").Value = "Unmeasureable"
$ws.Range("U40").Value = "Unmeasureable"

# Trial 5 (row 41)
$ws.Range("P41").Value = 1
$ws.Range("Q41").Value = 2
$ws.Range("R41").Value = 116
$ws.Range("S41").Value = 12393
$ws.Range("T41").Value = "Unmeasureable"
$ws.Range("U41").Value = "Unmeasureable"

# Update the view's top-left cell and selection to match
$ws.Activate()
$excel.Goto($ws.Range("I19"), $true)
$ws.Range("V38").Select()
